$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldTexts = @(
    "Supports the GCS if coverage is **Low** `nOther members: Global South + EU `n(25-33% of world emissions)",
    "Supports the GCS if coverage is **Mid** `nGlobal South + China `n(56% of world emissions)",
    "Supports the GCS if coverage is **High** `nGlobal South + China + EU + various HICs `n(UK, Japan, Korea, Canada...; 64-72% of emissions)",
    "Supports the GCS if coverage is **High**, **color** variant `nGlobal South + China + EU + various HICs `n+ Distributive effects shown using colors on world map"
)

$newTexts = @(
    "Supports the GCS if coverage is **Low**<br>Other members: Global South + EU<br>(25-33% of world emissions)",
    "Supports the GCS if coverage is **Mid**<br>Global South + China<br>(56% of world emissions)",
    "Supports the GCS if coverage is **High**<br>Global South + China + EU + various HICs<br>(UK, Japan, Korea, Canada...; 64-72% of emissions)",
    "Supports the GCS if coverage is **High**, **color** variant<br>Global South + China + EU + various HICs<br>+ Distributive effects shown using colors on world map"
)

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null) {
        for ($i = 0; $i -lt $oldTexts.Length; $i++) {
            if ($val -eq $oldTexts[$i]) {
                $cell.Value2 = $newTexts[$i]
                break
            }
        }
    }
}
